$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused 3rd data row (row 3) and the 7th column (G)
$ws.Rows.Item(3).Delete()
$ws.Columns.Item(7).Delete()

# Update header row (row 1)
$ws.Range("A1").Value = "bro"
$ws.Range("B1").Value = "totalClients"
$ws.Range("C1").Value = "assetsUnderCustody"
$ws.Range("D1").Value = "currentMarketValue"
$ws.Range("E1").Value = "unrealisedProfitLoss"
$ws.Range("F1").Value = "totalLedgerBalance"

# Update data row (row 2) - consolidated/aggregated summary row
$ws.Range("A2").Value = "N/A "
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 34219570.14
$ws.Range("D2").Value = 32371096.8
$ws.Range("E2").Value = -2097612.33
$ws.Range("F2").Value = 517639.77
